$wb = $excel.ActiveWorkbook

# Remove the "vol min", "vol max", "Expected part" and "Left/right" columns
# (columns C:F) from each check-protocol structures sheet. The "Mandatory"
# column (previously column G) slides left into column C.
foreach ($name in @("Clinical Structures", "opt structures", "couch_structures")) {
    $ws = $wb.Worksheets.Item($name)
    $rng = $ws.Range("C1:F1").EntireColumn
    [void]$rng.Select()
    $rng.Delete()
}

# The active sheet moves from "Clinical Structures" to "couch_structures".
$ws4 = $wb.Worksheets.Item("couch_structures")
$ws4.Activate()
